$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: fill in the "Work Done" and "Bugs / Limitations" notes for Day 7 ---
$ws.Range("C11").Value = "Spider now accepts input from a csv file that contains a list of companies to be crawled. Modified spider to get stockname from title of company page. Then used nsepy library to get stock info."
$ws.Range("E11").Value = "Company names need to be a perfect match. Improper company names will be skipped."
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12: new worklog entry for Day 8 (2019-06-26) ---
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 43642
$ws.Range("C12").Value = "built input.csv file to take in company name and keywords as input and store the output of keyword count in front of company name"
$ws.Range("E12").Value = "company names and keywords both need to be a perfect match. Improper company names will be skipped and mismatch keywords wont be counted"
$ws.Rows.Item(12).RowHeight = 90

# Match the date format / style already used by the other "Date" cells
$ws.Range("B12").NumberFormat = $ws.Range("B11").NumberFormat

# --- Update the view state to match where the author ended up after editing ---
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("E13").Select()
